$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.075.05"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.677.23"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "215.54"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +1.94%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "21.30"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +5.53%  "
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "1.912.85"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.675.31"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("E15").Value = "  +1.55%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "66.01"
$cell.Style = "Normal"
$ws.Range("D17").Value = "27.067.96"
$ws.Range("E17").Value = "  +0.45%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "237.37"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.38%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "8.14"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("E24").Value = "  -1.92%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "146.72"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.57%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "7.22"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("E27").Value = "  +2.45%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "1.548.79"
$ws.Range("E33").Value = "  +6.05%  "
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("E36").Value = "  +3.51%  "
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  +0.10%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "67.66"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.70%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.61"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -3.23%  "
$ws.Range("E44").Value = "  -2.06%  "
$ws.Range("D45").Value = "1.822.45"
$ws.Range("E45").Value = "  +0.71%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.782"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.12%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "90.73"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("E50").Value = "  +2.73%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "8.02"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +4.72%  "
